# Sync attendance_reports: in the "Recorded By" column (G) of the
# "Session Analysis Results" sheet, swap the order of "System" and the
# recorder's e-mail address, i.e. turn
#   "System, dnasr281@gmail.com"
# into
#   "dnasr281@gmail.com, System"
# wherever that exact text currently appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$col = 7   # column G = "Recorded By"
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
